# Update System, removed inverter
# This script updates the PV_Size, Battery_Size, PCM_Heating_Size, PCM_Cooling_Size,
# and the associated Training/Testing cost columns across worksheets Fold_1..Fold_5
# to reflect the new system configuration (inverter removed), per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 1664.473699999999
$ws.Range("C2").Value = 1990.3115
$ws.Range("D2").Value = 12.20105
$ws.Range("E2").Value = 775.7470999999998
$ws.Range("F2").Value = 993.8527999999997
$ws.Range("G2").Value = 993.8527999999997
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 41349.70080000001
$ws.Range("L2").Value = 993.8529999999998
$ws.Range("M2").Value = 40355.8478
$ws.Range("N2").Value = 2027.873
$ws.Range("O2").Value = 38327.9748
$ws.Range("B3").Value = 7985.395
$ws.Range("C3").Value = 2961.323
$ws.Range("E3").Value = 1652.713
$ws.Range("F3").Value = 3265.302
$ws.Range("G3").Value = 3120.128
$ws.Range("H3").Value = 145.175
$ws.Range("I3").Value = 145.175
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3483.02
$ws.Range("L3").Value = 3120.128
$ws.Range("M3").Value = 362.892
$ws.Range("N3").Value = 362.892
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 10243.474
$ws.Range("C4").Value = 3427.315
$ws.Range("E4").Value = 631.981
$ws.Range("F4").Value = 3649.922
$ws.Range("G4").Value = 3649.922101369863
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3864.975
$ws.Range("L4").Value = 3864.975
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 1994.3703
$ws.Range("C2").Value = 2056.10125
$ws.Range("E2").Value = 767.7413999999999
$ws.Range("F2").Value = 1105.34615
$ws.Range("G2").Value = 1105.34615
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 20860.6868
$ws.Range("L2").Value = 1105.346
$ws.Range("M2").Value = 19755.3408
$ws.Range("N2").Value = 1059.1138
$ws.Range("O2").Value = 18696.2272
$ws.Range("B3").Value = 6717.371
$ws.Range("C3").Value = 2938.913
$ws.Range("E3").Value = 1196.855
$ws.Range("F3").Value = 2983.484
$ws.Range("G3").Value = 2706.009
$ws.Range("H3").Value = 277.475
$ws.Range("I3").Value = 277.475
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6586.832399999999
$ws.Range("L3").Value = 2706.009
$ws.Range("M3").Value = 3880.8234
$ws.Range("N3").Value = 626.818
$ws.Range("O3").Value = 3254.0054
$ws.Range("B4").Value = 8628.054
$ws.Range("C4").Value = 3073.43
$ws.Range("E4").Value = 1250.284
$ws.Range("F4").Value = 3112.591
$ws.Range("G4").Value = 3112.590561643835
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3755.8482
$ws.Range("L4").Value = 3327.643
$ws.Range("M4").Value = 428.2052
$ws.Range("N4").Value = 428.2052
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 2157.95765
$ws.Range("C2").Value = 2000.9977
$ws.Range("D2").Value = 11.65735
$ws.Range("E2").Value = 789.1447499999999
$ws.Range("F2").Value = 1151.8421
$ws.Range("G2").Value = 1151.8421
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 17387.7698
$ws.Range("L2").Value = 1151.842
$ws.Range("M2").Value = 16235.9278
$ws.Range("N2").Value = 1560.6412
$ws.Range("O2").Value = 14675.287
$ws.Range("B3").Value = 7985.912
$ws.Range("C3").Value = 2921.255
$ws.Range("E3").Value = 1734.917
$ws.Range("F3").Value = 3355.068
$ws.Range("G3").Value = 3117.744
$ws.Range("H3").Value = 237.324
$ws.Range("I3").Value = 237.324
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3136.499800000001
$ws.Range("L3").Value = 3117.744
$ws.Range("M3").Value = 18.7558
$ws.Range("N3").Value = 18.7558
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 10243.474
$ws.Range("C4").Value = 3427.315
$ws.Range("E4").Value = 631.981
$ws.Range("F4").Value = 3649.922
$ws.Range("G4").Value = 3649.922101369863
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3864.975
$ws.Range("L4").Value = 3864.975
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 2451.1313
$ws.Range("C2").Value = 2121.3981
$ws.Range("D2").Value = 12.88205
$ws.Range("E2").Value = 751.8251
$ws.Range("F2").Value = 1256.86665
$ws.Range("G2").Value = 1256.86665
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3028.7838
$ws.Range("L2").Value = 1256.867
$ws.Range("M2").Value = 1771.9168
$ws.Range("N2").Value = 414.9076
$ws.Range("O2").Value = 1357.0092
$ws.Range("B3").Value = 7985.395
$ws.Range("C3").Value = 2961.323
$ws.Range("E3").Value = 1652.713
$ws.Range("F3").Value = 3356.025
$ws.Range("G3").Value = 3120.128
$ws.Range("H3").Value = 235.898
$ws.Range("I3").Value = 235.898
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3120.128
$ws.Range("L3").Value = 3120.128
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 10243.474
$ws.Range("C4").Value = 3427.315
$ws.Range("E4").Value = 631.981
$ws.Range("F4").Value = 3649.922
$ws.Range("G4").Value = 3649.922101369863
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3864.975
$ws.Range("L4").Value = 3864.975
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 2541.111450000001
$ws.Range("C2").Value = 2379.98045
$ws.Range("D2").Value = 7.543699999999999
$ws.Range("E2").Value = 873.9492499999999
$ws.Range("F2").Value = 1315.1723
$ws.Range("G2").Value = 1315.1723
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1707.1182
$ws.Range("L2").Value = 1315.172
$ws.Range("M2").Value = 391.9462
$ws.Range("N2").Value = 56.09
$ws.Range("O2").Value = 335.8562
$ws.Range("B3").Value = 7985.395
$ws.Range("C3").Value = 2961.323
$ws.Range("E3").Value = 1652.713
$ws.Range("F3").Value = 3356.025
$ws.Range("G3").Value = 3120.128
$ws.Range("H3").Value = 235.898
$ws.Range("I3").Value = 235.898
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3120.128
$ws.Range("L3").Value = 3120.128
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("B4").Value = 10243.474
$ws.Range("C4").Value = 3427.315
$ws.Range("E4").Value = 631.981
$ws.Range("F4").Value = 3649.922
$ws.Range("G4").Value = 3649.922101369863
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3864.975
$ws.Range("L4").Value = 3864.975
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0

